$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 735.1177
$ws.Range("I53").Value = 104.5
$ws.Range("K53").Value = 104.5
$ws.Range("M53").Value = 532.5

# Row 62
$ws.Range("H62").Value = 7860.273
$ws.Range("I62").Value = 7896.3
$ws.Range("K62").Value = 7896.3
$ws.Range("M62").Value = -7272.3

# Row 65
$ws.Range("H65").Value = 7860.273
$ws.Range("I65").Value = 7896.3
$ws.Range("K65").Value = 39481.5
$ws.Range("M65").Value = -36361.5

# Row 80
$ws.Range("H80").Value = 1146.28
$ws.Range("I80").Value = 899
$ws.Range("J80").Value = 1340.5714
$ws.Range("K80").Value = 2697
$ws.Range("L80").Value = 4021.7142
$ws.Range("M80").Value = -1699
$ws.Range("N80").Value = -6017.7142

# Row 83
$ws.Range("H83").Value = 1146.28
$ws.Range("I83").Value = 899
$ws.Range("J83").Value = 1340.5714
$ws.Range("K83").Value = 8091
$ws.Range("L83").Value = 12065.1426
$ws.Range("M83").Value = -3099
$ws.Range("N83").Value = -22049.1426

# Row 86
$ws.Range("H86").Value = 2542.8
$ws.Range("I86").Value = 1469
$ws.Range("J86").Value = 3079.7
$ws.Range("K86").Value = 1469
$ws.Range("L86").Value = 3079.7
$ws.Range("M86").Value = -346
$ws.Range("N86").Value = -5325.7

# Row 88
$ws.Range("H88").Value = 2068.0588
$ws.Range("I88").Value = 1371.25
$ws.Range("J88").Value = 2282.4614
$ws.Range("K88").Value = 1371.25
$ws.Range("L88").Value = 2282.4614
$ws.Range("M88").Value = -965.25
$ws.Range("N88").Value = -3094.4614

# Row 89
$ws.Range("H89").Value = 2542.8
$ws.Range("I89").Value = 1469
$ws.Range("J89").Value = 3079.7
$ws.Range("K89").Value = 7345
$ws.Range("L89").Value = 15398.5
$ws.Range("M89").Value = -1729
$ws.Range("N89").Value = -26630.5

# Row 91
$ws.Range("H91").Value = 2068.0588
$ws.Range("I91").Value = 1371.25
$ws.Range("J91").Value = 2282.4614
$ws.Range("K91").Value = 1371.25
$ws.Range("L91").Value = 2282.4614
$ws.Range("M91").Value = 32.75
$ws.Range("N91").Value = -5090.4614

# Row 106
$ws.Range("H106").Value = 3884.2
$ws.Range("J106").Value = 3980.75
$ws.Range("L106").Value = 3980.75
$ws.Range("N106").Value = -5242.75

# Row 136
$ws.Range("H136").Value = 35782.61
$ws.Range("J136").Value = 35782.61
$ws.Range("L136").Value = 35782.61
$ws.Range("N136").Value = -45982.61


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 52
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# Row 97
$ws.Range("H97").Value = 520.85187
$ws.Range("I97").Value = 535.7917
$ws.Range("J97").Value = 401.33334
$ws.Range("K97").Value = 535.7917
$ws.Range("L97").Value = 401.33334
$ws.Range("M97").Value = -39.79169999999999
$ws.Range("N97").Value = -1393.33334

# Row 117
$ws.Range("H117").Value = 400033100
$ws.Range("J117").Value = 400033100
$ws.Range("L117").Value = 400033100
$ws.Range("N117").Value = -400042278

# Row 119
$ws.Range("H119").Value = 53999.668
$ws.Range("J119").Value = 53999.668
$ws.Range("L119").Value = 53999.668
$ws.Range("N119").Value = -63675.668

# Row 120
$ws.Range("H120").Value = 36999
$ws.Range("J120").Value = 36999
$ws.Range("L120").Value = 36999
$ws.Range("N120").Value = -46675

# Row 121
$ws.Range("H121").Value = 89994.5
$ws.Range("J121").Value = 89994.5
$ws.Range("L121").Value = 89994.5
$ws.Range("N121").Value = -93488.5

# Row 122
$ws.Range("H122").Value = 3010.6853
$ws.Range("I122").Value = 2540.1777
$ws.Range("K122").Value = 7620.533100000001
$ws.Range("M122").Value = -5170.533100000001


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 7279
$ws.Range("I11").Value = 6500
$ws.Range("J11").Value = 8447.5
$ws.Range("K11").Value = 6500
$ws.Range("L11").Value = 8447.5
$ws.Range("M11").Value = -6360
$ws.Range("N11").Value = -8727.5

# Row 107
$ws.Range("H107").Value = 3838.375
$ws.Range("I107").Value = 3274.25
$ws.Range("J107").Value = 4402.5
$ws.Range("K107").Value = 3274.25
$ws.Range("L107").Value = 4402.5
$ws.Range("M107").Value = -1354.25
$ws.Range("N107").Value = -8242.5


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2199.6287
$ws.Range("I31").Value = 1596.3871
$ws.Range("J31").Value = 6874.75
$ws.Range("K31").Value = 1596.3871
$ws.Range("L31").Value = 6874.75
$ws.Range("M31").Value = -1301.3871
$ws.Range("N31").Value = -7464.75

# Row 34
$ws.Range("H34").Value = 2199.6287
$ws.Range("I34").Value = 1596.3871
$ws.Range("J34").Value = 6874.75
$ws.Range("K34").Value = 1596.3871
$ws.Range("L34").Value = 6874.75
$ws.Range("M34").Value = -1394.3871
$ws.Range("N34").Value = -7278.75


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 1180.1666
$ws.Range("I122").Value = 1196.5
$ws.Range("J122").Value = 1176.9
$ws.Range("K122").Value = 10768.5
$ws.Range("L122").Value = 10592.1
$ws.Range("M122").Value = -8318.5
$ws.Range("N122").Value = -15492.1

# Row 129
$ws.Range("H129").Value = 1043.5
$ws.Range("I129").Value = 511.0909
$ws.Range("K129").Value = 1533.2727
$ws.Range("M129").Value = 3466.7273

# Row 131
$ws.Range("H131").Value = 358461.6
$ws.Range("I131").Value = 1000800.7
$ws.Range("J131").Value = 1606.5555
$ws.Range("K131").Value = 3002402.1
$ws.Range("L131").Value = 4819.666499999999
$ws.Range("M131").Value = -2997362.1
$ws.Range("N131").Value = -14899.6665


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 94551.8
$ws.Range("J39").Value = 94551.8
$ws.Range("L39").Value = 94551.8
$ws.Range("N39").Value = -95615.8


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2554.389
$ws.Range("J22").Value = 3081
$ws.Range("L22").Value = 3081
$ws.Range("N22").Value = -3671

# Row 27
$ws.Range("H27").Value = 2554.389
$ws.Range("J27").Value = 3081
$ws.Range("L27").Value = 3081
$ws.Range("N27").Value = -3295

# Row 40
$ws.Range("H40").Value = 4406.7856
$ws.Range("I40").Value = 3963.3635
$ws.Range("K40").Value = 3963.3635
$ws.Range("M40").Value = -3827.3635

# Row 50
$ws.Range("H50").Value = 41907.125
$ws.Range("I50").Value = 40061
$ws.Range("J50").Value = 42522.5
$ws.Range("K50").Value = 40061
$ws.Range("L50").Value = 42522.5
$ws.Range("M50").Value = -39424
$ws.Range("N50").Value = -43796.5

# Row 55
$ws.Range("H55").Value = 413.7857
$ws.Range("I55").Value = 314.14285
$ws.Range("J55").Value = 513.4286
$ws.Range("K55").Value = 314.14285
$ws.Range("L55").Value = 513.4286
$ws.Range("M55").Value = -141.14285
$ws.Range("N55").Value = -859.4286

# Row 82
$ws.Range("H82").Value = 22666.8
$ws.Range("J82").Value = 3110.6667
$ws.Range("L82").Value = 3110.6667
$ws.Range("N82").Value = -3832.6667

# Row 85
$ws.Range("H85").Value = 22666.8
$ws.Range("J85").Value = 3110.6667
$ws.Range("L85").Value = 3110.6667
$ws.Range("N85").Value = -5606.6667


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2291.8076
$ws.Range("I132").Value = 2272.9565
$ws.Range("K132").Value = 6818.869499999999
$ws.Range("M132").Value = -4288.869499999999

# Row 135
$ws.Range("H135").Value = 70166.664
$ws.Range("J135").Value = 70166.664
$ws.Range("L135").Value = 70166.664
$ws.Range("N135").Value = -80306.664

# Row 136
$ws.Range("H136").Value = 9496.040000000001
$ws.Range("I136").Value = 9766.708000000001
$ws.Range("K136").Value = 29300.124
$ws.Range("M136").Value = -26750.124

